# Auto-generated Excel COM-interop script
# Applies updated profit/price figures to the Shinryu_Profits workbook sheets
# (covers 27 leve rows across 8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1599.1
$ws.Range("I32").Value = 1350
$ws.Range("K32").Value = 1350
$ws.Range("M32").Value = -1024

$ws.Range("H62").Value = 3342.2856
$ws.Range("J62").Value = 3117.7273
$ws.Range("L62").Value = 3117.7273
$ws.Range("N62").Value = -4365.7273

$ws.Range("H65").Value = 3342.2856
$ws.Range("J65").Value = 3117.7273
$ws.Range("L65").Value = 15588.6365
$ws.Range("N65").Value = -21828.6365

$ws.Range("H125").Value = 1854.7858
$ws.Range("I125").Value = 603.5
$ws.Range("J125").Value = 2793.25
$ws.Range("K125").Value = 5431.5
$ws.Range("L125").Value = 25139.25
$ws.Range("M125").Value = -2971.5
$ws.Range("N125").Value = -30059.25

$ws.Range("H137").Value = 8098880.5
$ws.Range("I137").Value = 14707700
$ws.Range("J137").Value = 73885.28999999999
$ws.Range("K137").Value = 44123100
$ws.Range("L137").Value = 221655.87
$ws.Range("M137").Value = -44120550
$ws.Range("N137").Value = -226755.87

$ws.Range("H138").Value = 3398.2188
$ws.Range("I138").Value = 2442.2144
$ws.Range("J138").Value = 3665.9
$ws.Range("K138").Value = 7326.6432
$ws.Range("L138").Value = 10997.7
$ws.Range("M138").Value = -2186.6432
$ws.Range("N138").Value = -21277.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1426.4286
$ws.Range("I2").Value = 1436.1538
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 1436.1538
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = -1323.1538
$ws.Range("N2").Value = -1526

$ws.Range("H61").Value = 2177.5862
$ws.Range("I61").Value = 2047.9166
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 2047.9166
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -1835.9166
$ws.Range("N61").Value = -3224

$ws.Range("H116").Value = 1426.4286
$ws.Range("I116").Value = 1436.1538
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 1436.1538
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = 857.8462
$ws.Range("N116").Value = -5888

$ws.Range("H136").Value = 2177.5862
$ws.Range("I136").Value = 2047.9166
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 6143.7498
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -3593.7498
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1426.4286
$ws.Range("I3").Value = 1436.1538
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 1436.1538
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = -1322.1538
$ws.Range("N3").Value = -1528

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2007.3167
$ws.Range("I31").Value = 1434.7273
$ws.Range("J31").Value = 2338.8157
$ws.Range("K31").Value = 1434.7273
$ws.Range("L31").Value = 2338.8157
$ws.Range("M31").Value = -1139.7273
$ws.Range("N31").Value = -2928.8157

$ws.Range("H34").Value = 2007.3167
$ws.Range("I34").Value = 1434.7273
$ws.Range("J34").Value = 2338.8157
$ws.Range("K34").Value = 1434.7273
$ws.Range("L34").Value = 2338.8157
$ws.Range("M34").Value = -1232.7273
$ws.Range("N34").Value = -2742.8157

$ws.Range("H58").Value = 2529.6553
$ws.Range("I58").Value = 1399.909
$ws.Range("J58").Value = 3220.0557
$ws.Range("K58").Value = 1399.909
$ws.Range("L58").Value = 3220.0557
$ws.Range("M58").Value = -1196.909
$ws.Range("N58").Value = -3626.0557

$ws.Range("H68").Value = 20035
$ws.Range("J68").Value = 20035
$ws.Range("L68").Value = 20035
$ws.Range("N68").Value = -21533

$ws.Range("H71").Value = 20035
$ws.Range("J71").Value = 20035
$ws.Range("L71").Value = 60105
$ws.Range("N71").Value = -67593

$ws.Range("H136").Value = 2529.6553
$ws.Range("I136").Value = 1399.909
$ws.Range("J136").Value = 3220.0557
$ws.Range("K136").Value = 4199.727000000001
$ws.Range("L136").Value = 9660.167099999999
$ws.Range("M136").Value = -1649.727000000001
$ws.Range("N136").Value = -14760.1671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1583.8551
$ws.Range("I68").Value = 1365.6
$ws.Range("J68").Value = 1707.8636
$ws.Range("K68").Value = 4096.799999999999
$ws.Range("L68").Value = 5123.5908
$ws.Range("M68").Value = -3285.799999999999
$ws.Range("N68").Value = -6745.5908

$ws.Range("H71").Value = 1583.8551
$ws.Range("I71").Value = 1365.6
$ws.Range("J71").Value = 1707.8636
$ws.Range("K71").Value = 12290.4
$ws.Range("L71").Value = 15370.7724
$ws.Range("M71").Value = -8234.4
$ws.Range("N71").Value = -23482.7724

$ws.Range("H122").Value = 8791.115
$ws.Range("I122").Value = 16712.309
$ws.Range("J122").Value = 869.9231
$ws.Range("K122").Value = 150410.781
$ws.Range("L122").Value = 7829.3079
$ws.Range("M122").Value = -147960.781
$ws.Range("N122").Value = -12729.3079

$ws.Range("H131").Value = 779.3200000000001
$ws.Range("I131").Value = 368.18182
$ws.Range("J131").Value = 1102.3572
$ws.Range("K131").Value = 1104.54546
$ws.Range("L131").Value = 3307.0716
$ws.Range("M131").Value = 3935.45454
$ws.Range("N131").Value = -13387.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6074.375
$ws.Range("I113").Value = 1741.3889
$ws.Range("J113").Value = 19073.334
$ws.Range("K113").Value = 1741.3889
$ws.Range("L113").Value = 19073.334
$ws.Range("M113").Value = 428.6111000000001
$ws.Range("N113").Value = -23413.334

$ws.Range("H132").Value = 2378.42
$ws.Range("I132").Value = 2283.225
$ws.Range("J132").Value = 2759.2
$ws.Range("K132").Value = 6849.674999999999
$ws.Range("L132").Value = 8277.599999999999
$ws.Range("M132").Value = -4319.674999999999
$ws.Range("N132").Value = -13337.6

$ws.Range("H133").Value = 36980
$ws.Range("J133").Value = 36980
$ws.Range("L133").Value = 36980
$ws.Range("N133").Value = -47100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5681.421
$ws.Range("I132").Value = 4380.6924
$ws.Range("J132").Value = 8499.666999999999
$ws.Range("K132").Value = 13142.0772
$ws.Range("L132").Value = 25499.001
$ws.Range("M132").Value = -10612.0772
$ws.Range("N132").Value = -30559.001

$ws.Range("H136").Value = 1293.0968
$ws.Range("I136").Value = 847.6667
$ws.Range("J136").Value = 4299.75
$ws.Range("K136").Value = 2543.0001
$ws.Range("L136").Value = 12899.25
$ws.Range("M136").Value = 6.999899999999798
$ws.Range("N136").Value = -17999.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

Write-Output "Applied 27 row updates across 8 sheets"
